# Insert a new price-record row (row 70) into the weekly "Ciruela" sheet,
# pushing the existing rows 70-111 down to 71-112.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("70").Insert()

$ws.Range("A70").Value = 11
$ws.Range("B70").Value = "Vega Monumental Concepción"
$ws.Range("C70").Value = "Bíobío"
$ws.Range("D70").Value = 44981
$ws.Range("E70").Value = 8
$ws.Range("F70").Value = "Fruta"
$ws.Range("G70").Value = 100103
$ws.Range("H70").Value = "Frutos de hueso (carozo)"
$ws.Range("I70").Value = 100103002
$ws.Range("J70").Value = "Ciruela"
$ws.Range("K70").Value = "Black Amber"
$ws.Range("L70").Value = "Primera"
$ws.Range("M70").Value = 220
$ws.Range("N70").Value = 9000
$ws.Range("O70").Value = 10000
$ws.Range("P70").Value = 9545
$ws.Range("Q70").Value = "`$/caja 18 kilos granel"
$ws.Range("R70").Value = "Provincia de Curicó"
$ws.Range("S70").Value = 530
$ws.Range("T70").Value = 18
